$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('A2').NumberFormat = "@"
$ws.Range('A2').Value = '2003611123'
$ws.Range('A2').Style = "Normal"
$ws.Range('B2').Value = 'Poco Смартфон C75 EU 6/128 ГБ, черный'
$ws.Range('C2').Value = 'Poco'
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '8947'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '9406'
$ws.Range('E2').Style = "Normal"
$ws.Range('F2').NumberFormat = "@"
$ws.Range('F2').Value = '39999'
$ws.Range('F2').Style = "Normal"
$ws.Range('H2').Value = '1 310 отзывов'
$ws.Range('I2').Value = 'SV SMARTPHONE'
$ws.Range('J2').Value = 'https://www.ozon.ru/seller/1685463/'
$ws.Range('K2').Value = 'ИП Соболь Влада Владимировна; 324290000003171'
$ws.Range('L2').Value = 'https://www.ozon.ru/product/poco-smartfon-c75-eu-6-128-gb-chernyy-2003611123/?at=gpt41jwomszPg4nQs589RxOsWQvnKEhDgkjBKcGB1DVz'

# Row 3
$ws.Range('A3').NumberFormat = "@"
$ws.Range('A3').Value = '2006925323'
$ws.Range('A3').Style = "Normal"
$ws.Range('B3').Value = 'Tecno Смартфон CAMON 40 Ростест (EAC) 8/256 ГБ, черный'
$ws.Range('C3').Value = 'Tecno'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '18314'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '18880'
$ws.Range('E3').Style = "Normal"
$ws.Range('F3').NumberFormat = "@"
$ws.Range('F3').Value = '23990'
$ws.Range('F3').Style = "Normal"
$ws.Range('G3').NumberFormat = "@"
$ws.Range('G3').Value = '4.9'
$ws.Range('G3').Style = "Normal"
$ws.Range('H3').Value = '208 отзывов'
$ws.Range('I3').Value = ""
$ws.Range('J3').Value = 'https://www.ozon.ru/seller/236588/'
$ws.Range('K3').Value = ""
$ws.Range('L3').Value = 'https://www.ozon.ru/product/tecno-smartfon-camon-40-rostest-eac-8-256-gb-chernyy-2006925323/?at=MZtvyLkwRfqOVN7jcgWw3BjtVXJ9vzCqDmJKPs2BwDvY'

# Row 4
$ws.Range('A4').NumberFormat = "@"
$ws.Range('A4').Value = '2133543492'
$ws.Range('A4').Style = "Normal"
$ws.Range('B4').Value = 'iQOO Смартфон Z10 Ростест (EAC) 8/256 ГБ, черный'
$ws.Range('C4').Value = 'iQOO'
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '20106'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '21219'
$ws.Range('E4').Style = "Normal"
$ws.Range('F4').NumberFormat = "@"
$ws.Range('F4').Value = '35999'
$ws.Range('F4').Style = "Normal"
$ws.Range('G4').NumberFormat = "@"
$ws.Range('G4').Value = '5'
$ws.Range('G4').Style = "Normal"
$ws.Range('H4').Value = '10 отзывов'
$ws.Range('I4').Value = 'iQOO Official Store'
$ws.Range('J4').Value = 'https://www.ozon.ru/seller/2898162/'
$ws.Range('K4').Value = 'ООО ООО "БАЙТ БУМ"127287, Россия, Москва, г Москва, ул Хуторская 2-я, стр 23, д 38А,; 131257700003538'
$ws.Range('L4').Value = 'https://www.ozon.ru/product/iqoo-smartfon-z10-rostest-eac-8-256-gb-chernyy-2133543492/?at=08tYX9g7McP6qr4QSEpQ9REc7kR46DsLwr80VIA2ownp'

# Row 5
$ws.Range('A5').NumberFormat = "@"
$ws.Range('A5').Value = '1771387609'
$ws.Range('A5').Style = "Normal"
$ws.Range('B5').Value = 'realme Смартфон Note 60X Ростест (EAC) 3/64 ГБ, черный'
$ws.Range('C5').Value = 'realme'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '4273'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '4999'
$ws.Range('E5').Style = "Normal"
$ws.Range('F5').NumberFormat = "@"
$ws.Range('F5').Value = '8499'
$ws.Range('F5').Style = "Normal"
$ws.Range('H5').Value = '7 034 отзыва'
$ws.Range('I5').Value = 'realme'
$ws.Range('J5').Value = 'https://www.ozon.ru/seller/216243/'
$ws.Range('K5').Value = 'ООО ООО "РМ КОММЬЮНИКЕЙШН"115280, г. Москва, вн.тер.г. Муниципальный округ Даниловский, улица Ленинская Слобода, Д. 19, помещение 21В/1H/5.1; 197746267355'
$ws.Range('L5').Value = 'https://www.ozon.ru/product/realme-smartfon-note-60x-rostest-eac-3-64-gb-chernyy-1771387609/?at=mqtko7PgAcEl1Y10coOx5McRM6BQJUojZmZ0UK6GZWY'

# Row 6
$ws.Range('A6').NumberFormat = "@"
$ws.Range('A6').Value = '1869776225'
$ws.Range('A6').Style = "Normal"
$ws.Range('B6').Value = 'Redmi Смартфон 8/256 ГБ, черный'
$ws.Range('C6').Value = 'Redmi'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '5901'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '6217'
$ws.Range('E6').Style = "Normal"
$ws.Range('F6').NumberFormat = "@"
$ws.Range('F6').Value = '13999'
$ws.Range('F6').Style = "Normal"
$ws.Range('G6').NumberFormat = "@"
$ws.Range('G6').Value = '4.4'
$ws.Range('G6').Style = "Normal"
$ws.Range('H6').Value = '466 отзывов'
$ws.Range('I6').Value = 'Bao Phone'
$ws.Range('J6').Value = 'https://www.ozon.ru/seller/2597173/'
$ws.Range('K6').Value = 'ИП Идалов Ибрагим Усманович; 321200000013285'
$ws.Range('L6').Value = 'https://www.ozon.ru/product/redmi-smartfon-8-256-gb-chernyy-1869776225/?at=DqtDLWJEnuj1wKBks45Lo6otB73YGlhPY3ZnoSlmOxMm'

# Row 7
$ws.Range('A7').NumberFormat = "@"
$ws.Range('A7').Value = '1711714454'
$ws.Range('A7').Style = "Normal"
$ws.Range('B7').Value = 'Tecno Смартфон Spark 30 5G Ростест (EAC) 6/128 ГБ, черный'
$ws.Range('C7').Value = 'Tecno'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '10011'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '10524'
$ws.Range('E7').Style = "Normal"
$ws.Range('F7').NumberFormat = "@"
$ws.Range('F7').Value = '19990'
$ws.Range('F7').Style = "Normal"
$ws.Range('G7').NumberFormat = "@"
$ws.Range('G7').Value = '4.9'
$ws.Range('G7').Style = "Normal"
$ws.Range('H7').Value = '971 отзыв'
$ws.Range('I7').Value = 'TECNO STORE'
$ws.Range('J7').Value = 'https://www.ozon.ru/seller/1452131/'
$ws.Range('K7').Value = 'ООО ООО "КИБЕР ЭНЕРДЖИ"127238, Россия, Москва, г Москва, Дмитровское шоссе, к 2, д; 321237700690864'
$ws.Range('L7').Value = 'https://www.ozon.ru/product/tecno-smartfon-spark-30-5g-rostest-eac-6-128-gb-chernyy-1711714454/?at=vQtrnPLY2tPD6WAphBNjgR5uQzJvkOu1oy4DZHpYNjy6'

# Row 8
$ws.Range('A8').NumberFormat = "@"
$ws.Range('A8').Value = '1897984574'
$ws.Range('A8').Style = "Normal"
$ws.Range('B8').Value = 'Смартфон 6,5-дюймовый противоударный смартфон Note60x, смартфон на базе Android 13, смартфон с большим объемом памяти, смартфон для студентов, поддержка русского языка Global 8/256 ГБ, черно-серый'
$ws.Range('C8').Value = 'Смартфоны'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3037'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '3233'
$ws.Range('E8').Style = "Normal"
$ws.Range('F8').NumberFormat = "@"
$ws.Range('F8').Value = '3773'
$ws.Range('F8').Style = "Normal"
$ws.Range('G8').NumberFormat = "@"
$ws.Range('G8').Value = '3.4'
$ws.Range('G8').Style = "Normal"
$ws.Range('H8').Value = '216 отзывов'
$ws.Range('I8').Value = 'Флагманский магазин UP'
$ws.Range('J8').Value = 'https://www.ozon.ru/seller/2489768/'
$ws.Range('K8').Value = 'Yaofengbaihuo'
$ws.Range('L8').Value = 'https://www.ozon.ru/product/smartfon-6-5-dyuymovyy-protivoudarnyy-smartfon-note60x-smartfon-na-baze-android-13-smartfon-1897984574/?at=GRt2NPOD0ck9ZDXzuNYNgNgHqy77nqiORKk4GIm79JBr'

# Row 9
$ws.Range('A9').NumberFormat = "@"
$ws.Range('A9').Value = '1743461395'
$ws.Range('A9').Style = "Normal"
$ws.Range('B9').Value = 'Tecno Смартфон Spark 30 Ростест (EAC) 8/256 ГБ, черный'
$ws.Range('C9').Value = 'Tecno'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '10633'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '11179'
$ws.Range('E9').Style = "Normal"
$ws.Range('F9').NumberFormat = "@"
$ws.Range('F9').Value = '18990'
$ws.Range('F9').Style = "Normal"
$ws.Range('H9').Value = '3 053 отзыва'
$ws.Range('I9').Value = 'TECNO STORE'
$ws.Range('J9').Value = 'https://www.ozon.ru/seller/1452131/'
$ws.Range('K9').Value = 'ООО ООО "КИБЕР ЭНЕРДЖИ"127238, Россия, Москва, г Москва, Дмитровское шоссе, к 2, д; 321237700690864'
$ws.Range('L9').Value = 'https://www.ozon.ru/product/tecno-smartfon-spark-30-rostest-eac-8-256-gb-chernyy-1743461395/?at=28t024ZBRf1L62mrCrRoX7nUl183plT4l4g3yfzPQ6rP'

# Row 10
$ws.Range('A10').NumberFormat = "@"
$ws.Range('A10').Value = '2137918430'
$ws.Range('A10').Style = "Normal"
$ws.Range('B10').Value = 'iQOO Смартфон Neo 10 Ростест (EAC) 16/512 ГБ, черный, черно-серый'
$ws.Range('C10').Value = 'iQOO'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '37482'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '39635'
$ws.Range('E10').Style = "Normal"
$ws.Range('F10').NumberFormat = "@"
$ws.Range('F10').Value = '52999'
$ws.Range('F10').Style = "Normal"
$ws.Range('H10').Value = '11 отзывов'
$ws.Range('I10').Value = 'iQOO Official Store'
$ws.Range('J10').Value = 'https://www.ozon.ru/seller/2898162/'
$ws.Range('K10').Value = 'ООО ООО "БАЙТ БУМ"127287, Россия, Москва, г Москва, ул Хуторская 2-я, стр 23, д 38А,; 131257700003538'
$ws.Range('L10').Value = 'https://www.ozon.ru/product/iqoo-smartfon-neo-10-rostest-eac-16-512-gb-chernyy-cherno-seryy-2137918430/?at=DqtDLWJEnu2BPmlPhOOk5kqckvYk78h9NmxwVckyxK85'

# Row 11
$ws.Range('A11').NumberFormat = "@"
$ws.Range('A11').Value = '1469526277'
$ws.Range('A11').Style = "Normal"
$ws.Range('B11').Value = 'Tecno Смартфон POVA 6 Pro 5G Ростест (EAC) 12/256 ГБ, черный'
$ws.Range('C11').Value = 'Tecno'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '22806'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '24069'
$ws.Range('E11').Style = "Normal"
$ws.Range('F11').NumberFormat = "@"
$ws.Range('F11').Value = '38337'
$ws.Range('F11').Style = "Normal"
$ws.Range('H11').Value = '2 221 отзыв'
$ws.Range('I11').Value = 'TECNO STORE'
$ws.Range('J11').Value = 'https://www.ozon.ru/seller/1452131/'
$ws.Range('K11').Value = 'ООО ООО "КИБЕР ЭНЕРДЖИ"127238, Россия, Москва, г Москва, Дмитровское шоссе, к 2, д; 321237700690864'
$ws.Range('L11').Value = 'https://www.ozon.ru/product/tecno-smartfon-pova-6-pro-5g-rostest-eac-12-256-gb-chernyy-1469526277/?at=Y7tjWvpnNClrrYVXuXvOzD3hYX1vywtNlAkYKTQo0Zql'

# Column width changes
$ws.Columns.Item(2).ColumnWidth = 197.15
$ws.Columns.Item(3).ColumnWidth = 10.15
$ws.Columns.Item(9).ColumnWidth = 23.15
$ws.Columns.Item(11).ColumnWidth = 155.15
$ws.Columns.Item(12).ColumnWidth = 180.15
